$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (moved one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the prices in column D for rows 31-33
$ws.Range("D31").Value = 3985.166
$ws.Range("D32").Value = 1646.582
$ws.Range("D33").Value = 3379.108

# Stash the formatting of B31:C33 first - re-merging cells below can
# normalize their border styling, so we restore it afterwards.
$ws.Range("B31:C33").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Re-apply the merged ranges that need to move so the stored mergeCells
# order matches the target order. (A1:D1, B30:C30 and A9:D9 are already
# in the right relative position so they are left untouched.)
$ws.Range("B33:C33").UnMerge()
$ws.Range("B32:C32").UnMerge()
$ws.Range("B31:C31").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A10:D10").UnMerge()

$ws.Range("B33:C33").Merge()
$ws.Range("B32:C32").Merge()
$ws.Range("B31:C31").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("A10:D10").Merge()

# Restore the original formatting for B31:C33
$ws.Range("H1:I3").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("H1:I3").Clear()
